$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3759.3142
$ws.Range("I80").Value = 408.86667
$ws.Range("J80").Value = 6272.15
$ws.Range("K80").Value = 1226.60001
$ws.Range("L80").Value = 18816.45
$ws.Range("M80").Value = -228.6000100000001
$ws.Range("N80").Value = -20812.45

$ws.Range("H83").Value = 3759.3142
$ws.Range("I83").Value = 408.86667
$ws.Range("J83").Value = 6272.15
$ws.Range("K83").Value = 3679.80003
$ws.Range("L83").Value = 56449.35
$ws.Range("M83").Value = 1312.19997
$ws.Range("N83").Value = -66433.35000000001

$ws.Range("H88").Value = 1580
$ws.Range("I88").Value = 1133.3334
$ws.Range("J88").Value = 2250
$ws.Range("K88").Value = 1133.3334
$ws.Range("L88").Value = 2250
$ws.Range("M88").Value = -727.3334
$ws.Range("N88").Value = -3062

$ws.Range("H91").Value = 1580
$ws.Range("I91").Value = 1133.3334
$ws.Range("J91").Value = 2250
$ws.Range("K91").Value = 1133.3334
$ws.Range("L91").Value = 2250
$ws.Range("M91").Value = 270.6666
$ws.Range("N91").Value = -5058

$ws.Range("H113").Value = 2856.125
$ws.Range("I113").Value = 2527.6
$ws.Range("J113").Value = 3005.4546
$ws.Range("K113").Value = 2527.6
$ws.Range("L113").Value = 3005.4546
$ws.Range("M113").Value = 726.4000000000001
$ws.Range("N113").Value = -9513.454600000001

$ws.Range("H132").Value = 1542.3016
$ws.Range("J132").Value = 2799.75
$ws.Range("L132").Value = 8399.25
$ws.Range("N132").Value = -13459.25

$ws.Range("H138").Value = 3654.2957
$ws.Range("J138").Value = 4254.518
$ws.Range("L138").Value = 12763.554
$ws.Range("N138").Value = -23043.554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 742
$ws.Range("I97").Value = 742
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 742
$ws.Range("M97").Value = -246
$ws.Range("N97").ClearContents()

$ws.Range("H122").Value = 4809873.5
$ws.Range("I122").Value = 3603
$ws.Range("K122").Value = 10809
$ws.Range("M122").Value = -8359

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2288.5806
$ws.Range("I86").Value = 2288.5806
$ws.Range("K86").Value = 2288.5806
$ws.Range("M86").Value = -1165.5806

$ws.Range("H89").Value = 2288.5806
$ws.Range("I89").Value = 2288.5806
$ws.Range("K89").Value = 11442.903
$ws.Range("M89").Value = -5826.902999999998

$ws.Range("H107").Value = 2564.65
$ws.Range("I107").Value = 2590.8333
$ws.Range("J107").Value = 2525.375
$ws.Range("K107").Value = 2590.8333
$ws.Range("L107").Value = 2525.375
$ws.Range("M107").Value = -670.8332999999998
$ws.Range("N107").Value = -6365.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 9800
$ws.Range("J11").Value = 9800
$ws.Range("L11").Value = 9800
$ws.Range("N11").Value = -10080

$ws.Range("H31").Value = 2902
$ws.Range("I31").Value = 2277.2593
$ws.Range("K31").Value = 2277.2593
$ws.Range("M31").Value = -1982.2593

$ws.Range("H34").Value = 2902
$ws.Range("I34").Value = 2277.2593
$ws.Range("K34").Value = 2277.2593
$ws.Range("M34").Value = -2075.2593

$ws.Range("H58").Value = 1979359.8
$ws.Range("I58").Value = 3032586.5
$ws.Range("J58").Value = 4559.375
$ws.Range("K58").Value = 3032586.5
$ws.Range("L58").Value = 4559.375
$ws.Range("M58").Value = -3032383.5
$ws.Range("N58").Value = -4965.375

$ws.Range("H117").Value = 60175
$ws.Range("J117").Value = 63900
$ws.Range("L117").Value = 63900
$ws.Range("N117").Value = -73078

$ws.Range("H136").Value = 1979359.8
$ws.Range("I136").Value = 3032586.5
$ws.Range("J136").Value = 4559.375
$ws.Range("K136").Value = 9097759.5
$ws.Range("L136").Value = 13678.125
$ws.Range("M136").Value = -9095209.5
$ws.Range("N136").Value = -18778.125

$ws.Range("H141").Value = 62892.2
$ws.Range("I141").Value = 27648
$ws.Range("J141").Value = 71703.25
$ws.Range("K141").Value = 27648
$ws.Range("L141").Value = 71703.25
$ws.Range("M141").Value = -22468
$ws.Range("N141").Value = -82063.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 94818.09
$ws.Range("I56").Value = 94818.09
$ws.Range("K56").Value = 94818.09
$ws.Range("M56").Value = -94288.09

$ws.Range("H68").Value = 7847.5
$ws.Range("J68").Value = 15187.714
$ws.Range("L68").Value = 45563.142
$ws.Range("N68").Value = -47185.142

$ws.Range("H71").Value = 7847.5
$ws.Range("J71").Value = 15187.714
$ws.Range("L71").Value = 136689.426
$ws.Range("N71").Value = -144801.426

$ws.Range("H137").Value = 35754236
$ws.Range("I137").Value = 50001930
$ws.Range("K137").Value = 150005790
$ws.Range("M137").Value = -150000690

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 4456.6
$ws.Range("I3").Value = 4301.5
$ws.Range("J3").Value = 4560
$ws.Range("K3").Value = 4301.5
$ws.Range("L3").Value = 4560
$ws.Range("M3").Value = -4185.5
$ws.Range("N3").Value = -4792

$ws.Range("H103").Value = 50000
$ws.Range("J103").Value = 50000
$ws.Range("L103").Value = 50000
$ws.Range("N103").Value = -52344

$ws.Range("H122").Value = 8635.875
$ws.Range("I122").Value = 18035.666
$ws.Range("J122").Value = 2996
$ws.Range("K122").Value = 54106.99800000001
$ws.Range("L122").Value = 8988
$ws.Range("M122").Value = -51656.99800000001
$ws.Range("N122").Value = -13888

$ws.Range("H132").Value = 5987.846
$ws.Range("I132").Value = 8457
$ws.Range("J132").Value = 3107.1667
$ws.Range("K132").Value = 25371
$ws.Range("L132").Value = 9321.500100000001
$ws.Range("M132").Value = -22841
$ws.Range("N132").Value = -14381.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2634.3333
$ws.Range("I93").Value = 2634.3333
$ws.Range("K93").Value = 2634.3333
$ws.Range("M93").Value = -1386.3333

$ws.Range("H99").Value = 53000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 53000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 53000
$ws.Range("N99").Value = -58990
$ws.Range("M99").ClearContents()

$ws.Range("H122").Value = 5326.85
$ws.Range("I122").Value = 4732.5
$ws.Range("K122").Value = 14197.5
$ws.Range("M122").Value = -11747.5

$ws.Range("H136").Value = 7338.231
$ws.Range("I136").Value = 6399.4
$ws.Range("J136").Value = 7925
$ws.Range("K136").Value = 19198.2
$ws.Range("L136").Value = 23775
$ws.Range("M136").Value = -16648.2
$ws.Range("N136").Value = -28875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H107").Value = 1016.76
$ws.Range("I107").Value = 422.46667
$ws.Range("K107").Value = 1267.40001
$ws.Range("M107").Value = 652.5999899999999

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H117").Value = 60000
$ws.Range("J117").Value = 60000
$ws.Range("L117").Value = 60000
$ws.Range("N117").Value = -69178

$ws.Range("H118").Value = 63000
$ws.Range("J118").Value = 63000
$ws.Range("L118").Value = 63000
$ws.Range("N118").Value = -66314

$ws.Range("H119").Value = 69690
$ws.Range("J119").Value = 69690
$ws.Range("L119").Value = 69690
$ws.Range("N119").Value = -79366

$ws.Range("H120").Value = 79400
$ws.Range("J120").Value = 79400
$ws.Range("L120").Value = 79400
$ws.Range("N120").Value = -89076

$ws.Range("H121").Value = 30420
$ws.Range("J121").Value = 30420
$ws.Range("L121").Value = 30420
$ws.Range("N121").Value = -33914

$ws.Range("H122").Value = 2666.158
$ws.Range("I122").Value = 2456
$ws.Range("J122").Value = 4452.5
$ws.Range("K122").Value = 7368
$ws.Range("L122").Value = 13357.5
$ws.Range("M122").Value = -4918
$ws.Range("N122").Value = -18257.5
